# Updated searchAny function (backend): the generic "/service/search" with
# an "any" body parameter is replaced by a query-string based
# "/service/search/" endpoint with one optional query parameter per field.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 58: GET /service/search  ->  GET /service/search/
$ws.Range("B58").Value = "/service/search/"
# The single "any" body parameter is gone.
$ws.Range("C58").Value = ""

# Rows 59-65: what used to be "Body" parameters (column C) are now shown
# as query-string fragments appended to the Path column (column B).
$ws.Range("B59").Value = "/username=:username"
$ws.Range("C59").Value = ""

$ws.Range("B60").Value = "&serviceName=:serviceName"
$ws.Range("C60").Value = ""

$ws.Range("B61").Value = "&category=:category"
$ws.Range("C61").Value = ""

$ws.Range("B62").Value = "&priceMin=:priceMin"
$ws.Range("C62").Value = ""

$ws.Range("B63").Value = "&priceMax=:priceMax"
$ws.Range("C63").Value = ""

$ws.Range("B64").Value = "&location=:location"
$ws.Range("C64").Value = ""

$ws.Range("B65").Value = "&description=:description"
$ws.Range("C65").Value = ""
